$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value, taken from the refreshed crypto price feed.
$updates = @{
    "D2" = "329.71"
    "E2" = "1.13%"
    "G2" = "9"
    "D3" = "43.94"
    "E3" = "-1.39%"
    "G3" = "9"
    "D4" = "5.488"
    "E4" = "-1.47%"
    "G4" = "9"
    "G5" = "9"
    "D6" = "1.984"
    "E6" = "4.17%"
    "G6" = "9"
    "E7" = "-4.67%"
    "G7" = "9"
    "D8" = "0.9512"
    "E8" = "1.12%"
    "G8" = "9"
    "E9" = "-4.39%"
    "G9" = "9"
    "D10" = "0.1879"
    "E10" = "1.03%"
    "G10" = "9"
    "D11" = "10.65"
    "E11" = "27.43%"
    "G11" = "9"
    "D12" = "0.09945"
    "E12" = "-0.08%"
    "G12" = "9"
    "D13" = "0.04717"
    "E13" = "10.61%"
    "G13" = "9"
    "E14" = "-0.33%"
    "G14" = "9"
    "D15" = "0.001273"
    "E15" = "-0.89%"
    "G15" = "9"
    "D16" = "0.04084"
    "E16" = "-2.65%"
    "G16" = "9"
    "D17" = "0.005978"
    "E17" = "-0.07%"
    "G17" = "9"
    "D18" = "3.369"
    "E18" = "-6.17%"
    "G18" = "9"
    "D19" = "4.384"
    "E19" = "1.73%"
    "G19" = "9"
    "D20" = "0.3469"
    "E20" = "-0.96%"
    "G20" = "9"
    "E21" = "3.13%"
    "G21" = "9"
    "D22" = "0.2587"
    "E22" = "2.40%"
    "G22" = "9"
    "D23" = "0.001267"
    "E23" = "2.37%"
    "G23" = "9"
    "D24" = "0.004357"
    "E24" = "-2.61%"
    "G24" = "9"
    "D25" = "0.0001200"
    "E25" = "1.69%"
    "G25" = "9"
    "D26" = "0.0003745"
    "E26" = "-6.14%"
    "G26" = "9"
    "G27" = "9"
    "G28" = "9"
    "G29" = "9"
    "G30" = "9"
    "G31" = "9"
    "G32" = "9"
    "G33" = "9"
    "G34" = "9"
    "G35" = "9"
    "G36" = "9"
    "G37" = "9"
    "D38" = "0.02583"
    "E38" = "-1.63%"
    "G38" = "9"
    "D39" = "0.05684"
    "E39" = "4.47%"
    "G39" = "9"
    "D40" = "0.007556"
    "E40" = "-1.63%"
    "G40" = "9"
    "D41" = "0.1397"
    "E41" = "0.16%"
    "G41" = "9"
    "D42" = "0.007407"
    "E42" = "5.31%"
    "G42" = "9"
    "D43" = "0.002015"
    "E43" = "-0.43%"
    "G43" = "9"
    "D44" = "0.008324"
    "E44" = "-2.83%"
    "G44" = "9"
    "D45" = "0.00007142"
    "E45" = "-0.95%"
    "G45" = "9"
    "D46" = "0.00000000750"
    "E46" = "0.00%"
    "G46" = "9"
    "D47" = "0.003531"
    "E47" = "55.50%"
    "G47" = "9"
    "D48" = "0.003561"
    "E48" = "0.81%"
    "G48" = "9"
    "D49" = "0.00002100"
    "E49" = "0.00%"
    "G49" = "9"
    "D50" = "0.0002000"
    "E50" = "0.00%"
    "G50" = "9"
    "G51" = "9"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (and the trailing
    # zeros / percent signs they carry) are not coerced into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop back to the default style so no stray text-format style index
    # is left behind on the cell.
    $cell.Style = "Normal"
}
